$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 314, shifting existing rows 314-324 down to 315-325.
$ws.Rows(314).Insert()

# Populate the newly inserted row 314 with the new record.
$ws.Cells.Item(314, 1).Value2 = 4
$ws.Cells.Item(314, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(314, 3).Value2 = "Los Lagos"
$ws.Cells.Item(314, 4).Value2 = 44747
$ws.Cells.Item(314, 5).Value2 = 10
$ws.Cells.Item(314, 6).Value2 = 100112045
$ws.Cells.Item(314, 7).Value2 = "Zapallo"
$ws.Cells.Item(314, 8).Value2 = "Paine"
$ws.Cells.Item(314, 9).Value2 = "1a (guarda)"
$ws.Cells.Item(314, 10).Value2 = 1000
$ws.Cells.Item(314, 11).Value2 = 500
$ws.Cells.Item(314, 12).Value2 = 500
$ws.Cells.Item(314, 13).Value2 = 500
$ws.Cells.Item(314, 14).Value2 = "$/kilo (volumen en unidades)"
$ws.Cells.Item(314, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(314, 16).Value2 = 500
$ws.Cells.Item(314, 17).Value2 = 1
$ws.Cells.Item(314, 18).Value2 = "Hortaliza"
